$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns AQ (PILA_UTM_E) / AR (PILA_UTM_N) previously held the placeholder
# text "NA" for every PILA tree row (there was no dSide/UTM position
# recorded yet). Real surveyed UTM easting/northing coordinates are now
# available for each tree, so replace the placeholder text with the actual
# numeric values, row by row.
$coords = @(
    @(2, 241805.0, 4190723.0),
    @(3, 241801.0, 4190716.0),
    @(4, 241798.0, 4190707.0),
    @(5, 241808.0, 4190714.0),
    @(8, 241787.0, 4190727.0),
    @(9, 241787.0, 4190723.0),
    @(10, 241784.0, 4190725.0),
    @(11, 241751.0, 4190727.0),
    @(15, 241756.0, 4190703.0),
    @(16, 241754.0, 4190669.0),
    @(18, 241761.0, 4190664.0),
    @(21, 241590.0, 4190673.0),
    @(22, 241646.0, 4190665.0),
    @(23, 241654.0, 4190690.0),
    @(24, 241654.0, 4190691.0),
    @(25, 241654.0, 4190693.0),
    @(26, 241652.0, 4190697.0),
    @(27, 241655.0, 4190698.0),
    @(28, 241655.0, 4190698.0),
    @(29, 241656.0, 4190699.0),
    @(30, 241670.0, 4190697.0),
    @(31, 241670.0, 4190698.0),
    @(32, 241670.0, 4190699.0),
    @(33, 241669.0, 4190696.0),
    @(34, 241670.0, 4190696.0),
    @(35, 214670.0, 4190697.0),
    @(36, 241675.0, 4190695.0),
    @(37, 241683.0, 4190704.0),
    @(38, 241682.0, 4190703.0),
    @(39, 241682.0, 4190702.0),
    @(40, 241685.0, 4190701.0),
    @(41, 241685.0, 4190693.0),
    @(42, 241685.0, 4190693.0),
    @(43, 241686.0, 4190704.0),
    @(44, 214670.0, 4190697.0)
)

# Row 44 (tree 47) previously carried a special highlighted format (a
# distinct font + fill) on AQ44/AR44 because those cells held the text note
# "use point 47 same clump" instead of coordinates. Now that row gets real
# numeric coordinates too, so its two cells should look like every other
# AQ/AR data cell in the column. Re-use the plain formatting already used
# by the rest of the column (e.g. AQ2) via a format-only paste so the
# workbook doesn't end up with a redundant, now-unused style.
$fmtSource = $ws.Range("AQ2")
$fmtSource.Copy()
$ws.Range("AQ44:AR44").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

foreach ($entry in $coords) {
    $row = $entry[0]
    $utmE = $entry[1]
    $utmN = $entry[2]
    $ws.Range("AQ$row").Value = $utmE
    $ws.Range("AR$row").Value = $utmN
}

# The plot note for tree 47 (BJ44) used to be split across two entries
# ("use point 47 same clump" in AQ44/AR44, and "SD, too old to assess" in
# BJ44). Since the duplicated waypoint note is gone, fold its meaning into
# the plot_notes cell so the context ("same clump as tree 34") isn't lost.
$ws.Range("BJ44").Value = "SD, too old to assess, same clump as tree 34"
